$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de usuario")

# Update the text of D12: extend the user story to mention returning to the main site
$ws.Range("D12").Value = "Necesito poder cerrar la sesion de mi cuenta logeada y retornar a la web principal"

# Update the view state to match: select E12 and scroll back to top-left (A1)
$ws.Range("A1").Select()
$ws.Range("E12").Select()
